$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.781.82"
$ws.Range("E2").Value = "  +1.18%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.114.92"
$ws.Range("E3").Value = "  -0.11%  "

$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.41"
$ws.Range("E5").Value = "  -0.42%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.22"
$ws.Range("E6").Value = "  +1.79%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.107.32"
$ws.Range("E8").Value = "  +0.13%  "

$ws.Range("E9").Value = "  -0.50%  "

$ws.Range("E10").Value = "  +9.85%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.76"
$ws.Range("E11").Value = "  +0.40%  "

$ws.Range("E12").Value = "  -0.78%  "

$ws.Range("E13").Value = "  +2.87%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.17"
$ws.Range("E14").Value = "  +4.48%  "

$ws.Range("E15").Value = "  -0.89%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.630.39"
$ws.Range("E16").Value = "  -0.09%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.688.03"
$ws.Range("E17").Value = "  +1.11%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.13"
$ws.Range("E18").Value = "  -2.05%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.111.54"
$ws.Range("E19").Value = "  -0.05%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "464.79"
$ws.Range("E20").Value = "  +2.45%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.31"
$ws.Range("E21").Value = "  +1.48%  "

$ws.Range("E22").Value = "  -0.68%  "

$ws.Range("E23").Value = "  -0.45%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.12"
$ws.Range("E24").Value = "  -3.74%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.87"
$ws.Range("E25").Value = "  -0.25%  "

$ws.Range("E26").Value = "  -0.10%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.94"

$ws.Range("E28").Value = "  -0.62%  "

$ws.Range("E29").Value = "  -1.61%  "

$ws.Range("E30").Value = "  -0.07%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "26.95"
$ws.Range("E32").Value = "  -0.79%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0872"
$ws.Range("E34").Value = "  +8.00%  "

$ws.Range("E35").Value = "  +2.29%  "

$ws.Range("E36").Value = "  +0.70%  "

$ws.Range("E37").Value = "  +12.15%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.08"
$ws.Range("E38").Value = "  -0.04%  "

$ws.Range("E39").Value = "  -0.03%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "447.19"
$ws.Range("E40").Value = "  +4.39%  "

$ws.Range("E41").Value = "  -1.52%  "

$ws.Range("E42").Value = "  -1.05%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.877.36"
$ws.Range("E43").Value = "  -2.82%  "

$ws.Range("E44").Value = "  +0.40%  "

$ws.Range("E45").Value = "  -0.56%  "

$ws.Range("E46").Value = "  -0.35%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "35.80"
$ws.Range("E47").Value = "  +3.15%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "123.46"
$ws.Range("E49").Value = "  -1.50%  "

$ws.Range("E50").Value = "  -0.63%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.65"
$ws.Range("E51").Value = "  -0.65%  "
